$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.051.23"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "2.055.06"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.52"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.75%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0813"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "2.358.73"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.755"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "2.065.16"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").Value = "37.964.62"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.11%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0219"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "1.486.64"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("E43").Value = "  +3.61%  "
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "2.246.65"
$ws.Range("E51").Value = "  +1.48%  "
